# Generate Report for Handback
# Update the handoff/handback timestamp columns with freshly generated
# report times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
# (this equals the later of the two locales' handoff datetimes below).
$wsOverview.Range("G2").Value = "2016-11-14 06:52:06"

# zh-cn detail sheet, first file row: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-11-14 06:51:52"
$wsZhCn.Range("K2").Value = "2016-11-14 06:52:43"

# de-de detail sheet, first file row: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-11-14 06:52:06"
$wsDeDe.Range("K2").Value = "2016-11-14 06:53:02"
